# Adds more "rotations" of the name cycle below the existing data.
# The existing 5x5 block (A1:E5) is read out in row-major order, which
# reveals a repeating cycle of 14 names. We keep extending the grid with
# that same cycle, adding as many new rows as requested (5 here, to grow
# the range from A1:E5 to A1:E10), so the user can run as many rotations
# as they want in one go.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$numCols = 5
$existingRows = 5
$rowsToAdd = 5

# Read the existing values in row-major order to recover the rotation cycle.
$flat = New-Object System.Collections.ArrayList
for ($r = 1; $r -le $existingRows; $r++) {
    for ($c = 1; $c -le $numCols; $c++) {
        [void]$flat.Add($ws.Cells.Item($r, $c).Value2)
    }
}

# Detect the smallest repeating period p such that flat[i] == flat[i - p]
# for every i >= p (the existing data need not be an exact multiple of the
# period, e.g. 25 cells of a 14-name cycle).
$n = $flat.Count
$cycleLen = $n
for ($p = 1; $p -lt $n; $p++) {
    $isPeriod = $true
    for ($i = $p; $i -lt $n; $i++) {
        if ($flat[$i] -ne $flat[$i - $p]) {
            $isPeriod = $false
            break
        }
    }
    if ($isPeriod) {
        $cycleLen = $p
        break
    }
}

$totalRows = $existingRows + $rowsToAdd
$totalCells = $totalRows * $numCols

for ($i = $flat.Count; $i -lt $totalCells; $i++) {
    $row = [math]::Floor($i / $numCols) + 1
    $col = ($i % $numCols) + 1
    $value = $flat[$i % $cycleLen]
    $ws.Cells.Item($row, $col).Value2 = $value
    [void]$flat.Add($value)
}
